# Updates the "cryptos" price/volume table with refreshed scrape values.
# Many of the Price/Volume(1h) cells look numeric (e.g. "1.00", "0.522")
# but must stay plain text (matching the existing inline-string cells),
# so we force the cell to Text format before writing, then restore the
# default ("Normal") style so no stray number-format override is left
# behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '66.190.02'
Set-TextValue $ws 'E2' '  +0.24%  '
Set-TextValue $ws 'D3' '3.071.56'
Set-TextValue $ws 'E3' '  +3.49%  '
Set-TextValue $ws 'D4' '1.00'
Set-TextValue $ws 'D5' '578.43'
Set-TextValue $ws 'E5' '  +0.61%  '
Set-TextValue $ws 'D6' '167.12'
Set-TextValue $ws 'E6' '  +3.16%  '
Set-TextValue $ws 'E7' '  -0.12%  '
Set-TextValue $ws 'D8' '3.069.97'
Set-TextValue $ws 'E8' '  +3.81%  '
Set-TextValue $ws 'D9' '0.522'
Set-TextValue $ws 'E10' '  +0.57%  '
Set-TextValue $ws 'E11' '  +0.36%  '
Set-TextValue $ws 'D12' '0.480'
Set-TextValue $ws 'E12' '  +5.80%  '
Set-TextValue $ws 'D13' '0.0000249'
Set-TextValue $ws 'E13' '  +0.92%  '
Set-TextValue $ws 'D14' '36.73'
Set-TextValue $ws 'E14' '  +7.03%  '
Set-TextValue $ws 'E15' '  -0.23%  '
Set-TextValue $ws 'D16' '3.576.82'
Set-TextValue $ws 'E16' '  +3.23%  '
Set-TextValue $ws 'D17' '66.187.52'
Set-TextValue $ws 'E17' '  +0.10%  '
Set-TextValue $ws 'D18' '7.17'
Set-TextValue $ws 'E18' '  +4.14%  '
Set-TextValue $ws 'D19' '3.068.06'
Set-TextValue $ws 'E19' '  +3.14%  '
Set-TextValue $ws 'D20' '16.12'
Set-TextValue $ws 'E20' '  +17.36%  '
Set-TextValue $ws 'D21' '464.15'
Set-TextValue $ws 'E21' '  +3.63%  '
Set-TextValue $ws 'E22' '  +5.17%  '
Set-TextValue $ws 'D23' '7.43'
Set-TextValue $ws 'E23' '  +3.08%  '
Set-TextValue $ws 'D24' '83.05'
Set-TextValue $ws 'E24' '  +1.32%  '
Set-TextValue $ws 'D25' '12.80'
Set-TextValue $ws 'E25' '  +4.82%  '
Set-TextValue $ws 'E26' '  +2.06%  '
Set-TextValue $ws 'E27' '  +0.60%  '
Set-TextValue $ws 'E28' '  +0.08%  '
Set-TextValue $ws 'D29' '8.09'
Set-TextValue $ws 'E29' '  -0.39%  '
Set-TextValue $ws 'D30' '2.42'
Set-TextValue $ws 'E30' '  +0.50%  '
Set-TextValue $ws 'D31' '2.66'
Set-TextValue $ws 'E31' '  +2.83%  '
Set-TextValue $ws 'E32' '  +1.62%  '
Set-TextValue $ws 'D33' '28.24'
Set-TextValue $ws 'E33' '  +3.97%  '
Set-TextValue $ws 'E34' '  +5.15%  '
Set-TextValue $ws 'D35' '0.999'
Set-TextValue $ws 'E35' '  -0.02%  '
Set-TextValue $ws 'D36' '0.996'
Set-TextValue $ws 'E36' '  +1.47%  '
Set-TextValue $ws 'D37' '5.87'
Set-TextValue $ws 'E37' '  +2.45%  '
Set-TextValue $ws 'D38' '48.81'
Set-TextValue $ws 'E38' '  +11.17%  '
Set-TextValue $ws 'D39' '49.92'
Set-TextValue $ws 'E39' '  +1.07%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws 'D40' '0.313'
Set-TextValue $ws 'E40' '  +4.25%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws 'D41' '2.03'
Set-TextValue $ws 'E41' '  +0.38%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 'D42' '2.90'
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws 'D43' '0.122'
Set-TextValue $ws 'E43' '  +2.47%  '
Set-TextValue $ws 'D44' '8.67'
Set-TextValue $ws 'E44' '  +3.65%  '
Set-TextValue $ws 'D45' '0.0360'
Set-TextValue $ws 'E45' '  +1.66%  '
Set-TextValue $ws 'D46' '380.79'
Set-TextValue $ws 'E46' '  -1.96%  '
Set-TextValue $ws 'D47' '2.758.63'
Set-TextValue $ws 'E47' '  +1.48%  '
Set-TextValue $ws 'D48' '134.65'
Set-TextValue $ws 'E48' '  +3.05%  '
Set-TextValue $ws 'E49' '  +0.01%  '
Set-TextValue $ws 'D50' '24.44'
Set-TextValue $ws 'E50' '  +5.68%  '
Set-TextValue $ws 'E51' '  +4.30%  '
